$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 68
$ws.Range("A68").Value = "2012.8.9.5"
$ws.Range("B68").Value = "制作 开发者展示界面"
$ws.Range("C68").Value = "界面风格不统一"
$ws.Range("D68").Value = 4

# Row 69
$ws.Range("A69").Value = "2012.8.9.9"
$ws.Range("B69").Value = "小组会议"

# Row 70
$ws.Range("A70").Value = "2012.8.9.16"
$ws.Range("B70").Value = "修改loading界面 初步制作教学图片"
$ws.Range("C70").Value = "教学界面需要调整"
$ws.Range("D70").Value = 5

# Row 71
$ws.Range("A71").Value = "2012.8.9.20"
$ws.Range("B71").Value = "完成开始弹出教学界面;初步制作暂停时教学界面;调整商店等图标"
$ws.Range("C71").Value = "`n"
$ws.Range("D71").Value = 4

# Row 72
$ws.Range("A72").Value = "2012.8.9.22"
$ws.Range("B72").Value = "完成暂停时教学界面，调整暂停界面摆设"
$ws.Range("D72").Value = 5

# Row 73
$ws.Range("B73").Value = "调整奖励，出属性球策略"

# Row 74
$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = "2012.10.23"
$ws.Range("A74").ClearFormats()
$ws.Range("B74").Value = "添加熊猫"
$ws.Range("D74").Value = 4

# Row 75
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "2012.10.24"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = "修改纹理图"
$ws.Range("D75").Value = 3

# Row 76
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "2012.11.3"
$ws.Range("A76").ClearFormats()
$ws.Range("B76").Value = "整理z轴"
$ws.Range("D76").Value = 3

# Row 77
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = "2012.11.4"
$ws.Range("A77").ClearFormats()
$ws.Range("B77").Value = "添加仓库可用和不可用"
$ws.Range("D77").Value = 3

# Wrap text + row height for the blank/newline note cell
$ws.Range("C71").WrapText = $true
$ws.Rows.Item(71).RowHeight = 27

# Update sheet view (scroll position + selection)
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("D79").Select()
